$wb = $excel.ActiveWorkbook

# --- Rename the original (sole) sheet --------------------------------------
$wb.Worksheets.Item(1).Name = "sabit_kalemler"

# --- Add the three new sheets, then move each into its final tab position --
# (Worksheets.Add() always inserts at the front, so every new sheet has to
#  be moved right after the sheet that should precede it. Worksheet handles
#  returned earlier go stale once the tab order changes, so every lookup
#  below re-resolves the sheet by name via $wb.Worksheets.Item(...).)

$ws = $wb.Worksheets.Add()
$ws.Name = "kilavuzluk"
$wb.Worksheets.Item("kilavuzluk").Move($null, $wb.Worksheets.Item("sabit_kalemler"))

$ws = $wb.Worksheets.Add()
$ws.Name = "romorkor_istanbul"
$wb.Worksheets.Item("romorkor_istanbul").Move($null, $wb.Worksheets.Item("kilavuzluk"))

$ws = $wb.Worksheets.Add()
$ws.Name = "romorkor_canakkale"
$wb.Worksheets.Item("romorkor_canakkale").Move($null, $wb.Worksheets.Item("romorkor_istanbul"))

# --- Header row values -------------------------------------------------
$wb.Worksheets.Item("kilavuzluk").Range("A1").Value = "hizmet_turu"
$wb.Worksheets.Item("kilavuzluk").Range("B1").Value = "taban"
$wb.Worksheets.Item("kilavuzluk").Range("C1").Value = "ilave"

$wb.Worksheets.Item("romorkor_istanbul").Range("A1").Value = "alt_boy"
$wb.Worksheets.Item("romorkor_istanbul").Range("B1").Value = "ust_boy"
$wb.Worksheets.Item("romorkor_istanbul").Range("C1").Value = "cins"
$wb.Worksheets.Item("romorkor_istanbul").Range("D1").Value = "ucret"

$wb.Worksheets.Item("romorkor_canakkale").Range("A1").Value = "alt_boy"
$wb.Worksheets.Item("romorkor_canakkale").Range("B1").Value = "ust_boy"
$wb.Worksheets.Item("romorkor_canakkale").Range("C1").Value = "cins"
$wb.Worksheets.Item("romorkor_canakkale").Range("D1").Value = "ucret"

# --- Copy the existing bold/centered/bordered header style onto the new
#     header rows (reuses the workbook's existing header cell style rather
#     than fabricating new ones). -----------------------------------------
$wb.Worksheets.Item("sabit_kalemler").Range("A1:B1").Copy()
$wb.Worksheets.Item("kilavuzluk").Range("A1:C1").PasteSpecial(-4122)

$wb.Worksheets.Item("sabit_kalemler").Range("A1:B1").Copy()
$wb.Worksheets.Item("romorkor_istanbul").Range("A1:D1").PasteSpecial(-4122)

$wb.Worksheets.Item("sabit_kalemler").Range("A1:B1").Copy()
$wb.Worksheets.Item("romorkor_canakkale").Range("A1:D1").PasteSpecial(-4122)

$wb.Worksheets.Item("sabit_kalemler").Select()
